$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B18").Value = "Filter Solar"
$ws.Range("C18").Value = "kecil"
$ws.Range("D18").Value = "2 pc"
$ws.Range("E18").Value = 425325
$ws.Range("F18").Formula = "=E18+10000"
$ws.Range("G18").Value = "27/03/2024"

$ws.Range("G19").Select()
